$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

$wsOverview.Range("G3").Value = "2016-08-24 04:46:14"

$wsZhCn.Range("H3").Value = "2016-08-24 04:46:08"
$wsZhCn.Range("K3").Value = "2016-08-24 04:46:26"

$wsDeDe.Range("H3").Value = "2016-08-24 04:46:14"
$wsDeDe.Range("K3").Value = "2016-08-24 04:46:33"
